$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.325.12"
$ws.Range("D2").NumberFormat = "General"
$ws.Range("E2").Value = "  +0.50%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.423.60"
$ws.Range("D3").NumberFormat = "General"
$ws.Range("E3").Value = "  +0.09%  "

$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "414.13"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = "  +0.75%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "128.73"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = "  -0.59%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.626"
$ws.Range("D7").NumberFormat = "General"
$ws.Range("E7").Value = "  -2.42%  "

$ws.Range("E8").Value = "  +0.09%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.729"
$ws.Range("D9").NumberFormat = "General"
$ws.Range("E9").Value = "  -1.48%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.140"
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").Value = "  -1.14%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "42.79"
$ws.Range("D11").NumberFormat = "General"
$ws.Range("E11").Value = "  +0.02%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000219"
$ws.Range("D12").NumberFormat = "General"
$ws.Range("E12").Value = "  -0.19%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "9.22"
$ws.Range("D13").NumberFormat = "General"
$ws.Range("E13").Value = "  +0.69%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.966.84"
$ws.Range("D14").NumberFormat = "General"
$ws.Range("E14").Value = "  +0.35%  "

$ws.Range("E15").Value = "  -0.35%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "20.54"
$ws.Range("D16").NumberFormat = "General"
$ws.Range("E16").Value = "  -3.38%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.430.31"
$ws.Range("D17").NumberFormat = "General"
$ws.Range("E17").Value = "  -0.12%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.80"
$ws.Range("D18").NumberFormat = "General"
$ws.Range("E18").Value = "  +5.94%  "

$ws.Range("E19").Value = "  -0.97%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "62.315.08"
$ws.Range("D20").NumberFormat = "General"
$ws.Range("E20").Value = "  +0.55%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "477.65"
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").Value = "  +7.13%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "91.99"
$ws.Range("D22").NumberFormat = "General"
$ws.Range("E22").Value = "  +0.23%  "

$ws.Range("E23").Value = "  +2.85%  "

$ws.Range("E24").Value = "  -0.35%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.30"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").Value = "  +1.04%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.68"
$ws.Range("D26").NumberFormat = "General"
$ws.Range("E26").Value = "  +9.31%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "33.44"
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").Value = "  -0.26%  "

$ws.Range("E28").Value = "  +0.50%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.73"
$ws.Range("D29").NumberFormat = "General"
$ws.Range("E29").Value = "  +1.48%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.65"
$ws.Range("D30").NumberFormat = "General"
$ws.Range("E30").Value = "  -3.34%  "

$ws.Range("E31").Value = "  -1.20%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.166"
$ws.Range("D32").NumberFormat = "General"
$ws.Range("E32").Value = "  -1.62%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.111"
$ws.Range("D33").NumberFormat = "General"
$ws.Range("E33").Value = "  -3.50%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "41.02"
$ws.Range("D34").NumberFormat = "General"
$ws.Range("E34").Value = "  -4.28%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "58.11"
$ws.Range("D36").NumberFormat = "General"
$ws.Range("E36").Value = "  +8.19%  "

$ws.Range("E37").Value = "  -2.60%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.999"
$ws.Range("D38").NumberFormat = "General"
$ws.Range("E38").Value = "  +0.09%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.06"
$ws.Range("D39").NumberFormat = "General"
$ws.Range("E39").Value = "  +4.39%  "

$ws.Range("E40").Value = "  -0.23%  "

$ws.Range("E41").Value = "  +2.13%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "147.40"
$ws.Range("D42").NumberFormat = "General"
$ws.Range("E42").Value = "  +4.22%  "

$ws.Range("E43").Value = "  -1.84%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.67"
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").Value = "  +11.02%  "

$ws.Range("E45").Value = "  +4.51%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.30"
$ws.Range("D46").NumberFormat = "General"
$ws.Range("E46").Value = "  +1.34%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.34"
$ws.Range("D47").NumberFormat = "General"
$ws.Range("E47").Value = "  +18.52%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "16.33"
$ws.Range("D48").NumberFormat = "General"
$ws.Range("E48").Value = "  -1.81%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0₃0533"
$ws.Range("D49").NumberFormat = "General"
$ws.Range("E49").Value = "  +23.90%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "22.28"
$ws.Range("D50").NumberFormat = "General"
$ws.Range("E50").Value = "  -0.23%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "113.64"
$ws.Range("D51").NumberFormat = "General"
$ws.Range("E51").Value = "  +8.05%  "
